$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V (match + odds) details between the two rows in each pair. ---
# A-E (index, country, tournament, season, match date) stay tied to the row position.

# Swap row 13 <-> row 14
$ws.Range("F13").Value = 'Horsens'
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 'Koge'
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1.72
$ws.Range("K13").Value = '01/08/2023 05:42'
$ws.Range("L13").Value = 1.81
$ws.Range("M13").Value = '04/08/2023 18:52'
$ws.Range("N13").Value = 3.95
$ws.Range("O13").Value = '01/08/2023 05:42'
$ws.Range("P13").Value = 4.16
$ws.Range("Q13").Value = '04/08/2023 18:56'
$ws.Range("R13").Value = 4.12
$ws.Range("S13").Value = '01/08/2023 05:42'
$ws.Range("T13").Value = 3.96
$ws.Range("U13").Value = '04/08/2023 18:56'
$ws.Range("V13").Value = 'https://www.betexplorer.com/football/denmark/1st-division/horsens-koge/zDaoiXfI/'
$ws.Range("F14").Value = 'Hillerod'
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 'Hobro'
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 2.23
$ws.Range("K14").Value = '29/07/2023 14:12'
$ws.Range("L14").Value = 3.15
$ws.Range("M14").Value = '04/08/2023 18:23'
$ws.Range("N14").Value = 3.38
$ws.Range("O14").Value = '29/07/2023 14:12'
$ws.Range("P14").Value = 3.49
$ws.Range("Q14").Value = '04/08/2023 18:23'
$ws.Range("R14").Value = 3.23
$ws.Range("S14").Value = '29/07/2023 14:12'
$ws.Range("T14").Value = 2.27
$ws.Range("U14").Value = '04/08/2023 18:23'
$ws.Range("V14").Value = 'https://www.betexplorer.com/football/denmark/1st-division/hillerod-hobro/E3bkjiAO/'

# Swap row 43 <-> row 45
$ws.Range("F43").Value = 'Hillerod'
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 'Sonderjyske'
$ws.Range("I43").Value = 2
$ws.Range("J43").Value = 3.9
$ws.Range("K43").Value = '28/08/2023 18:42'
$ws.Range("L43").Value = 4.16
$ws.Range("M43").Value = '01/09/2023 18:58'
$ws.Range("N43").Value = 3.8
$ws.Range("O43").Value = '28/08/2023 18:42'
$ws.Range("P43").Value = 3.85
$ws.Range("Q43").Value = '01/09/2023 18:58'
$ws.Range("R43").Value = 1.79
$ws.Range("S43").Value = '28/08/2023 18:42'
$ws.Range("T43").Value = 1.83
$ws.Range("U43").Value = '01/09/2023 18:58'
$ws.Range("V43").Value = 'https://www.betexplorer.com/football/denmark/1st-division/hillerod-sonderjyske/EwHtnuEE/'
$ws.Range("F45").Value = 'Vendsyssel'
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 'Fredericia'
$ws.Range("I45").Value = 2
$ws.Range("J45").Value = 2.02
$ws.Range("K45").Value = '27/08/2023 13:12'
$ws.Range("L45").Value = 2.46
$ws.Range("M45").Value = '01/09/2023 18:55'
$ws.Range("N45").Value = 3.79
$ws.Range("O45").Value = '27/08/2023 13:12'
$ws.Range("P45").Value = 3.9
$ws.Range("Q45").Value = '01/09/2023 18:55'
$ws.Range("R45").Value = 3.19
$ws.Range("S45").Value = '27/08/2023 13:12'
$ws.Range("T45").Value = 2.62
$ws.Range("U45").Value = '01/09/2023 18:55'
$ws.Range("V45").Value = 'https://www.betexplorer.com/football/denmark/1st-division/vendsyssel-ff-fredericia/zmIxmLb8/'

# Swap row 68 <-> row 69
$ws.Range("F68").Value = 'B.93'
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 'Hillerod'
$ws.Range("I68").Value = 3
$ws.Range("J68").Value = 2.6
$ws.Range("K68").Value = '30/09/2023 13:12'
$ws.Range("L68").Value = 3.6
$ws.Range("M68").Value = '06/10/2023 18:59'
$ws.Range("N68").Value = 3.5
$ws.Range("O68").Value = '30/09/2023 13:12'
$ws.Range("P68").Value = 3.81
$ws.Range("Q68").Value = '06/10/2023 18:59'
$ws.Range("R68").Value = 2.5
$ws.Range("S68").Value = '30/09/2023 13:12'
$ws.Range("T68").Value = 1.98
$ws.Range("U68").Value = '06/10/2023 18:59'
$ws.Range("V68").Value = 'https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-hillerod/hCneYENl/'
$ws.Range("F69").Value = 'Horsens'
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 'Aalborg'
$ws.Range("I69").Value = 4
$ws.Range("J69").Value = 3.66
$ws.Range("K69").Value = '29/09/2023 18:13'
$ws.Range("L69").Value = 3.88
$ws.Range("M69").Value = '06/10/2023 18:53'
$ws.Range("N69").Value = 3.86
$ws.Range("O69").Value = '29/09/2023 18:13'
$ws.Range("P69").Value = 3.72
$ws.Range("Q69").Value = '06/10/2023 18:54'
$ws.Range("R69").Value = 1.84
$ws.Range("S69").Value = '29/09/2023 18:13'
$ws.Range("T69").Value = 1.93
$ws.Range("U69").Value = '06/10/2023 18:54'
$ws.Range("V69").Value = 'https://www.betexplorer.com/football/denmark/1st-division/horsens-aalborg/SUv8VC77/'

# Swap row 79 <-> row 80
$ws.Range("F79").Value = 'B.93'
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 'Sonderjyske'
$ws.Range("I79").Value = 4
$ws.Range("J79").Value = 5.01
$ws.Range("K79").Value = '22/10/2023 15:12'
$ws.Range("L79").Value = 8.029999999999999
$ws.Range("M79").Value = '27/10/2023 18:58'
$ws.Range("N79").Value = 4.47
$ws.Range("O79").Value = '22/10/2023 15:12'
$ws.Range("P79").Value = 5.46
$ws.Range("Q79").Value = '27/10/2023 18:58'
$ws.Range("R79").Value = 1.57
$ws.Range("S79").Value = '22/10/2023 15:12'
$ws.Range("T79").Value = 1.34
$ws.Range("U79").Value = '27/10/2023 18:50'
$ws.Range("V79").Value = 'https://www.betexplorer.com/football/denmark/1st-division/boldklubben-1893-sonderjyske/2PEpqWy7/'
$ws.Range("F80").Value = 'Hobro'
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 'Koge'
$ws.Range("I80").Value = 1
$ws.Range("J80").Value = 1.74
$ws.Range("K80").Value = '22/10/2023 16:12'
$ws.Range("L80").Value = 1.65
$ws.Range("M80").Value = '27/10/2023 18:51'
$ws.Range("N80").Value = 3.92
$ws.Range("O80").Value = '22/10/2023 16:12'
$ws.Range("P80").Value = 4.26
$ws.Range("Q80").Value = '27/10/2023 18:51'
$ws.Range("R80").Value = 4.42
$ws.Range("S80").Value = '22/10/2023 16:12'
$ws.Range("T80").Value = 4.85
$ws.Range("U80").Value = '27/10/2023 18:51'
$ws.Range("V80").Value = 'https://www.betexplorer.com/football/denmark/1st-division/hobro-koge/UBGxoh7f/'

# --- Append two new match rows (85, 86), copying formatting from the last existing row. ---
$ws.Range("A84:V84").Copy()
$ws.Range("A85:V86").PasteSpecial(-4122)

# Row 85
$ws.Range("A85").Value = 84
$ws.Range("B85").Value = 'denmark'
$ws.Range("C85").Value = '1st-division'
$ws.Range("D85").Value = '2023-2024'
$ws.Range("E85").Value = 45233.79166666666
$ws.Range("F85").Value = 'Koge'
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 'Horsens'
$ws.Range("I85").Value = 4
$ws.Range("J85").Value = 3.16
$ws.Range("K85").Value = '27/10/2023 20:13'
$ws.Range("L85").Value = 3.3
$ws.Range("M85").Value = '03/11/2023 18:58'
$ws.Range("N85").Value = 3.7
$ws.Range("O85").Value = '27/10/2023 20:13'
$ws.Range("P85").Value = 3.7
$ws.Range("Q85").Value = '03/11/2023 18:58'
$ws.Range("R85").Value = 2.14
$ws.Range("S85").Value = '27/10/2023 20:13'
$ws.Range("T85").Value = 2.12
$ws.Range("U85").Value = '03/11/2023 18:58'
$ws.Range("V85").Value = 'https://www.betexplorer.com/football/denmark/1st-division/koge-horsens/dMGHoPc5/'

# Row 86
$ws.Range("A86").Value = 85
$ws.Range("B86").Value = 'denmark'
$ws.Range("C86").Value = '1st-division'
$ws.Range("D86").Value = '2023-2024'
$ws.Range("E86").Value = 45233.79166666666
$ws.Range("F86").Value = 'Sonderjyske'
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 'Hobro'
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1.54
$ws.Range("K86").Value = '27/10/2023 20:13'
$ws.Range("L86").Value = 1.53
$ws.Range("M86").Value = '03/11/2023 18:54'
$ws.Range("N86").Value = 4.36
$ws.Range("O86").Value = '27/10/2023 20:13'
$ws.Range("P86").Value = 4.29
$ws.Range("Q86").Value = '03/11/2023 18:54'
$ws.Range("R86").Value = 5.47
$ws.Range("S86").Value = '27/10/2023 20:13'
$ws.Range("T86").Value = 6.04
$ws.Range("U86").Value = '03/11/2023 18:54'
$ws.Range("V86").Value = 'https://www.betexplorer.com/football/denmark/1st-division/sonderjyske-hobro/jqGLpqDB/'

